# "valuation of kalimasta give to user committee"
# Zero-out the quantity columns that had been counted for the RCC/rebar
# work items so the estimate reflects the portion handed over to the
# user committee (the underlying dimension formulas are kept, just
# multiplied by 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C30").Formula = "=0*4"
$ws.Range("C31").Formula = "=0*3"
$ws.Range("C32").Formula = "=0*1"

$ws.Range("C44").Formula = "=0*4"
$ws.Range("C45").Formula = "=0*3"
$ws.Range("C46").Formula = "=0*1"

$ws.Range("C59").Formula = "=0*2"
$ws.Range("C60").Formula = "=0*3"
$ws.Range("C61").Formula = "=0*2*TRUNC((D59-3*0.72)/0.125,0)"
$ws.Range("C62").Formula = "=0*3*TRUNC((D60-2*0.72)/0.125,0)"
$ws.Range("C63").Formula = "=0*TRUNC(12.5/0.5,0)+1"
$ws.Range("C64").Formula = "=0*TRUNC((25.333-0.333-0.75*3)/0.5,0)+1"
$ws.Range("C65").Formula = "=0*30"
$ws.Range("C66").Formula = "=0*9*2"
$ws.Range("C67").Formula = "=0*7*4"
$ws.Range("C68").Formula = "=0*16"
$ws.Range("C69").Formula = "=0*8"

# Move the saved cursor/view position to where the author last left it.
[void]$ws.Range("G73").Select()
